# Hortaliza, Femacal de La Calera - Apio: insert one new daily price
# record at row 367, pushing the existing rows (367-487) down by one
# (to 368-488), matching the "Fruta / hortaliza, semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 367; Excel shifts rows
# 367..487 down to 368..488 and copies formatting (incl. the date
# style on column D) from the row above, same as a manual row insert.
$ws.Rows.Item(367).Insert()

# Populate the newly inserted row 367 with the new record.
$ws.Range("A367").Value = 3
$ws.Range("B367").Value = "Femacal de La Calera"
$ws.Range("C367").Value = "Coquimbo"
$ws.Range("D367").Value = 44876
$ws.Range("E367").Value = 5
$ws.Range("F367").Value = 100112017
$ws.Range("G367").Value = "Apio"
$ws.Range("H367").Value = "Americana (o)"
$ws.Range("I367").Value = "Primera"
$ws.Range("J367").Value = 130
$ws.Range("K367").Value = 9000
$ws.Range("L367").Value = 9000
$ws.Range("M367").Value = 9000
$ws.Range("N367").Value = "`$/docena de matas"
$ws.Range("O367").Value = "Provincia de Limarí"
$ws.Range("P367").Value = 1500
$ws.Range("Q367").Value = 6
$ws.Range("R367").Value = "Hortaliza"
